$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 6.992417114397731
$ws.Range("C3").Value = 1.052572528399653
$ws.Range("E3").Value = 2.396905458966625
$ws.Range("C4").Value = 0.475547144596522
$ws.Range("E4").Value = 1.495774185788745
$ws.Range("C5").Value = 4.312104569761943
$ws.Range("E5").Value = 1.976172106438545
$ws.Range("C6").Value = 4.639864760432189
$ws.Range("E6").Value = 4.074582884048139
$ws.Range("C7").Value = 1.418316910291906
$ws.Range("E7").Value = 3.247860853607465
$ws.Range("C8").Value = 3.520945360626571
$ws.Range("E8").Value = 2.547371915279606
$ws.Range("C9").Value = 3.285568146716344
$ws.Range("E9").Value = 3.38738696315446
$ws.Range("C10").Value = 2.441258738366514
$ws.Range("E10").Value = 3.091878630346012
$ws.Range("C11").Value = 2.480871685520603
$ws.Range("E11").Value = 2.296583397191387
$ws.Range("C12").Value = 3.071095202329288
$ws.Range("E12").Value = 2.781797072072023
$ws.Range("C13").Value = 2.994116795316071
$ws.Range("E13").Value = 2.939737488252936
$ws.Range("C14").Value = 2.671604274379558
$ws.Range("E14").Value = 2.997455747043043
$ws.Range("C15").Value = 2.08524086077817
$ws.Range("E15").Value = 2.197771900625956
$ws.Range("C16").Value = 4.939003803830477
$ws.Range("E16").Value = 3.440178795466697
$ws.Range("C17").Value = 2.93530792557688
$ws.Range("E17").Value = 3.215749572764803
$ws.Range("C18").Value = 1.635353376270698
$ws.Range("E18").Value = 2.328770194687713
$ws.Range("C19").Value = 1.984020855913604
$ws.Range("E19").Value = 1.88544721086894
